$d = $word.ActiveDocument

# wdHeaderFooterIndex constants
$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2

foreach ($sec in $d.Sections) {
    foreach ($hfIndex in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage)) {

        # Headers: the BTec logo picture is renamed image1.jpg -> image2.jpg
        $hdr = $sec.Headers.Item($hfIndex)
        if ($hdr.Exists) {
            $ils = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $ils.Count; $i++) {
                $shp = $ils.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }

        # Footers: the Pearson Edexcel logo picture is renamed image2.png -> image1.png
        $ftr = $sec.Footers.Item($hfIndex)
        if ($ftr.Exists) {
            $ils = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $ils.Count; $i++) {
                $shp = $ils.Item($i)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
